$d = $word.ActiveDocument

function Merge-ParagraphRuns($paraIndex, $firstRunLength, $fullText) {
    $p = $d.Paragraphs($paraIndex)
    $pRange = $p.Range
    $start = $pRange.Start
    $end = $pRange.End - 1   # exclude the paragraph mark

    # Keep the paragraph's original first run (characters 0..firstRunLength)
    # untouched so its xml:space="preserve" run survives, delete the rest
    # of the paragraph's runs, then append the remaining text onto that
    # same (now sole) run via InsertAfter.
    $restRange = $d.Range($start + $firstRunLength, $end)
    $restRange.Delete()

    $remainder = $fullText.Substring($firstRunLength)
    $keepRange2 = $d.Range($start, $start + $firstRunLength)
    $keepRange2.InsertAfter($remainder)
}

Merge-ParagraphRuns 1 8 "Answers: Introduction to vectors"
Merge-ParagraphRuns 2 5 "Zheng Chen"
Merge-ParagraphRuns 4 7 "Answers to questions relating to the guide on introduction to vectors."
